# 自动更新Excel文件 - 2026-01-15 23:14:11
# For every data row, recompute the "剩余" (remaining) days (column E) based on
# the "总天" (total days, column D) and the "开始时间" (start date, column F),
# as of "today" = 2026-01-16. When the remaining days would drop to zero or
# below, the cycle is treated as renewed: the start date is reset to today and
# the remaining days reset back to the total.

function Get-DayNumber {
    param([int]$y, [int]$m, [int]$d)
    if ($m -le 2) {
        $y = $y - 1
        $mAdj = $m + 12
    } else {
        $mAdj = $m
    }
    $era = [Math]::Floor($y / 400)
    $yoe = $y - $era * 400
    $mp = ($mAdj - 3) % 12
    $doy = [Math]::Floor((153 * $mp + 2) / 5) + $d - 1
    $doe = $yoe * 365 + [Math]::Floor($yoe / 4) - [Math]::Floor($yoe / 100) + $doy
    return $era * 146097 + $doe
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$todayYear = 2026
$todayMonth = 1
$todayDay = 16
$todayNum = Get-DayNumber $todayYear $todayMonth $todayDay
$todayValue = $todayYear * 10000 + $todayMonth * 100 + $todayDay

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $dCell = $ws.Cells.Item($r, 4)
    $eCell = $ws.Cells.Item($r, 5)
    $fCell = $ws.Cells.Item($r, 6)

    $total = $dCell.Value2
    $startRaw = $fCell.Value2

    if ($total -eq $null -or $startRaw -eq $null) {
        continue
    }

    $startText = [string]([int]$startRaw)
    if ($startText.Length -ne 8) {
        # Not a well-formed yyyyMMdd date value - leave this row untouched.
        continue
    }

    $year = [int]$startText.Substring(0, 4)
    $month = [int]$startText.Substring(4, 2)
    $day = [int]$startText.Substring(6, 2)

    if ($month -lt 1 -or $month -gt 12 -or $day -lt 1 -or $day -gt 31) {
        continue
    }

    $startNum = Get-DayNumber $year $month $day
    $elapsed = $todayNum - $startNum
    $remaining = [int]$total - $elapsed

    if ($remaining -le 0) {
        # Cycle exhausted: restart it as of today.
        $eCell.Value2 = [int]$total
        $fCell.Value2 = [int]$todayValue
    } else {
        $eCell.Value2 = $remaining
    }
}
